$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44: add a new time-log entry (date + start time), matching
# the "integrated expo" entry that was previously blank.
$ws.Range("A44").Value = 43071
$ws.Range("B44").Value = 0.91666666666666663
